# Insert a new data row at row 106 (pushing existing rows 106..184 down to 107..185)
# and populate it with the new price-record values for Zapallo italiano
# at Vega Monumental Concepción (Bíobío).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 106-184 down by inserting a new blank row at 106.
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with values.
$ws.Range("A106").Value = 11
$ws.Range("B106").Value = "Vega Monumental Concepción"
$ws.Range("C106").Value = "Bíobío"
$ws.Range("D106").Value = 44944
$ws.Range("E106").Value = 8
$ws.Range("F106").Value = 100112032
$ws.Range("G106").Value = "Zapallo italiano"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 350
$ws.Range("K106").Value = 9000
$ws.Range("L106").Value = 10000
$ws.Range("M106").Value = 9429
$ws.Range("N106").Value = "$/caja 50 unidades"
$ws.Range("O106").Value = "Región Metropolitana"
$ws.Range("P106").Value = 189
$ws.Range("Q106").Value = 50
$ws.Range("R106").Value = "Hortaliza"
